$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row.
# All data rows (2 through 351) change from serial 45175 (2023-09-06)
# to serial 45177 (2023-09-08).
$oldSerial = 45175
$newSerial = 45177

$firstRow = 2
$lastRow = 351

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
